# Recompute the Bmp2 (ligand, sending-cluster) and Eng (receptor, target-cluster)
# average/total expression values with the refreshed TPM inputs, and write the
# derived specificity / edge-weight columns that follow from them.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.370286
$ws.Range("H2").Value = 7.110858
$ws.Range("I2").Value = 0.3026841782318013
$ws.Range("J2").Value = 0.3026841782318014
$ws.Range("M2").Value = 218.721583
$ws.Range("N2").Value = 656.164749
$ws.Range("O2:P2").Value = 0.7793342808141792
$ws.Range("Q2").Value = 518.432706082738
$ws.Range("R2").Value = 4665.894354744642
$ws.Range("S2").Value = 0.2358921563561117
$ws.Range("T2").Value = 0.2358921563561118

# Row 3
$ws.Range("G3").Value = 2.370286
$ws.Range("H3").Value = 7.110858
$ws.Range("I3").Value = 0.3026841782318013
$ws.Range("J3").Value = 0.3026841782318014
$ws.Range("O3:P3").Value = 0.164954193449581
$ws.Range("Q3").Value = 109.7316658525793
$ws.Range("R3").Value = 987.5849926732141
$ws.Range("S3").Value = 0.04992902449017601
$ws.Range("T3").Value = 0.04992902449017602

# Row 4
$ws.Range("G4").Value = 2.370286
$ws.Range("H4").Value = 7.110858
$ws.Range("I4").Value = 0.3026841782318013
$ws.Range("J4").Value = 0.3026841782318014
$ws.Range("M4").Value = 8.515309999999999
$ws.Range("N4").Value = 25.54593
$ws.Range("O4:P4").Value = 0.03034118948727519
$ws.Range("Q4").Value = 20.18372007866
$ws.Range("R4").Value = 181.65348070794
$ws.Range("S4").Value = 0.00918379800653126
$ws.Range("T4").Value = 0.009183798006531262

# Row 5
$ws.Range("G5").Value = 2.370286
$ws.Range("H5").Value = 7.110858
$ws.Range("I5").Value = 0.3026841782318013
$ws.Range("J5").Value = 0.3026841782318014
$ws.Range("M5").Value = 7.120231
$ws.Range("N5").Value = 21.360693
$ws.Range("O5:P5").Value = 0.02537033624896462
$ws.Range("Q5").Value = 16.876983856066
$ws.Range("R5").Value = 151.892854704594
$ws.Range("S5:T5").Value = 0.007679199378982338

# Row 6
$ws.Range("I6:J6").Value = 0.2022126055089961
$ws.Range("M6").Value = 218.721583
$ws.Range("N6").Value = 656.164749
$ws.Range("O6:P6").Value = 0.7793342808141792
$ws.Range("Q6").Value = 346.3465744740263
$ws.Range("R6").Value = 3117.119170266237
$ws.Range("S6:T6").Value = 0.1575912154859148

# Row 7
$ws.Range("I7:J7").Value = 0.2022126055089961
$ws.Range("O7:P7").Value = 0.164954193449581
$ws.Range("S7").Value = 0.03335581724707475
$ws.Range("T7").Value = 0.03335581724707476

# Row 8
$ws.Range("I8:J8").Value = 0.2022126055089961
$ws.Range("M8").Value = 8.515309999999999
$ws.Range("N8").Value = 25.54593
$ws.Range("O8:P8").Value = 0.03034118948727519
$ws.Range("Q8").Value = 13.48403028467667
$ws.Range("R8").Value = 121.35627256209
$ws.Range("S8:T8").Value = 0.006135370980464078

# Row 9
$ws.Range("I9:J9").Value = 0.2022126055089961
$ws.Range("M9").Value = 7.120231
$ws.Range("N9").Value = 21.360693
$ws.Range("O9:P9").Value = 0.02537033624896462
$ws.Range("Q9").Value = 11.27491664283433
$ws.Range("R9").Value = 101.474249785509
$ws.Range("S9:T9").Value = 0.005130201795542468

# Row 10
$ws.Range("G10").Value = 2.286703333333333
$ws.Range("H10").Value = 6.860109999999999
$ws.Range("I10:J10").Value = 0.2920107190904054
$ws.Range("M10").Value = 218.721583
$ws.Range("N10").Value = 656.164749
$ws.Range("O10:P10").Value = 0.7793342808141792
$ws.Range("Q10").Value = 500.1513729180433
$ws.Range("R10").Value = 4501.36235626239
$ws.Range("S10:T10").Value = 0.2275739637523524

# Row 11
$ws.Range("G11").Value = 2.286703333333333
$ws.Range("H11").Value = 6.860109999999999
$ws.Range("I11:J11").Value = 0.2920107190904054
$ws.Range("O11:P11").Value = 0.164954193449581
$ws.Range("Q11").Value = 105.8622318476811
$ws.Range("R11").Value = 952.7600866291299
$ws.Range("S11:T11").Value = 0.04816839264618999

# Row 12
$ws.Range("G12").Value = 2.286703333333333
$ws.Range("H12").Value = 6.860109999999999
$ws.Range("I12:J12").Value = 0.2920107190904054
$ws.Range("M12").Value = 8.515309999999999
$ws.Range("N12").Value = 25.54593
$ws.Range("O12:P12").Value = 0.03034118948727519
$ws.Range("Q12").Value = 19.47198776136666
$ws.Range("R12").Value = 175.2478898522999
$ws.Range("S12").Value = 0.008859952560237479
$ws.Range("T12").Value = 0.008859952560237477

# Row 13
$ws.Range("G13").Value = 2.286703333333333
$ws.Range("H13").Value = 6.860109999999999
$ws.Range("I13:J13").Value = 0.2920107190904054
$ws.Range("M13").Value = 7.120231
$ws.Range("N13").Value = 21.360693
$ws.Range("O13:P13").Value = 0.02537033624896462
$ws.Range("Q13").Value = 16.28185596180333
$ws.Range("R13").Value = 146.53670365623
$ws.Range("S13").Value = 0.007408410131625539
$ws.Range("T13").Value = 0.007408410131625538

# Row 14
$ws.Range("G14").Value = 1.590394666666667
$ws.Range("H14").Value = 4.771184
$ws.Range("I14:J14").Value = 0.2030924971687972
$ws.Range("M14").Value = 218.721583
$ws.Range("N14").Value = 656.164749
$ws.Range("O14:P14").Value = 0.7793342808141792
$ws.Range("Q14").Value = 347.8536390880907
$ws.Range("R14").Value = 3130.682751792816
$ws.Range("S14:T14").Value = 0.1582769452198003

# Row 15
$ws.Range("G15").Value = 1.590394666666667
$ws.Range("H15").Value = 4.771184
$ws.Range("I15:J15").Value = 0.2030924971687972
$ws.Range("O15:P15").Value = 0.164954193449581
$ws.Range("Q15").Value = 73.62683496269689
$ws.Range("R15").Value = 662.641514664272
$ws.Range("S15:T15").Value = 0.03350095906614026

# Row 16
$ws.Range("G16").Value = 1.590394666666667
$ws.Range("H16").Value = 4.771184
$ws.Range("I16:J16").Value = 0.2030924971687972
$ws.Range("M16").Value = 8.515309999999999
$ws.Range("N16").Value = 25.54593
$ws.Range("O16:P16").Value = 0.03034118948727519
$ws.Range("Q16").Value = 13.54270360901333
$ws.Range("R16").Value = 121.88433248112
$ws.Range("S16:T16").Value = 0.006162067940042375

# Row 17
$ws.Range("G17").Value = 1.590394666666667
$ws.Range("H17").Value = 4.771184
$ws.Range("I17:J17").Value = 0.2030924971687972
$ws.Range("M17").Value = 7.120231
$ws.Range("N17").Value = 21.360693
$ws.Range("O17:P17").Value = 0.02537033624896462
$ws.Range("Q17").Value = 11.32397740783467
$ws.Range("R17").Value = 101.915796670512
$ws.Range("S17").Value = 0.00515252494281428
$ws.Range("T17").Value = 0.005152524942814279
